$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2, "H 72"), shifting all subsequent rows up by one
$ws.Rows.Item(2).Delete()

# Apply the updated missing-value mask: clear newly-missing cells and fill newly-available ones
$ws.Range("C2").Value = 10
$ws.Range("C3").ClearContents()
$ws.Range("F3").Value = 0.70917
$ws.Range("D6").Value = -15.4
$ws.Range("F6").ClearContents()
$ws.Range("F7").Value = 0.71266
$ws.Range("D8").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("F15").Value = 0.70964
$ws.Range("D16").Value = -14.1
$ws.Range("D18").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("C20").Value = 13.5
$ws.Range("D20").Value = -15.3
$ws.Range("C21").ClearContents()
$ws.Range("C22").Value = 11.5
$ws.Range("D22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C24").Value = 12.5
$ws.Range("F24").Value = 0.7101499999999999
$ws.Range("C25").ClearContents()
$ws.Range("F28").Value = 0.70963
$ws.Range("D30").Value = -13.8
$ws.Range("F30").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("F32").Value = 0.71028
$ws.Range("F33").Value = 0.72961
$ws.Range("F34").ClearContents()
$ws.Range("F35").ClearContents()
$ws.Range("F36").ClearContents()
$ws.Range("F49").Value = 0.70924
$ws.Range("C52").Value = 10.8
$ws.Range("F52").ClearContents()
$ws.Range("C53").ClearContents()
$ws.Range("C56").Value = 11.9
$ws.Range("C57").ClearContents()
$ws.Range("C58").Value = 11.2
$ws.Range("C59").ClearContents()

Write-Host "done"
